$wb = $excel.ActiveWorkbook

# Updates to column F ("想去人数" / interest count) across sheets 1, 2, 3, 4
# Parallel arrays: sheet index, row number, expected old value, new value
$sheetIdx = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 2, 2, 3, 3, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4, 4)
$rowNum   = @(2, 3, 4, 5, 6, 8, 9, 10, 12, 13, 14, 15, 17, 20, 22, 6, 7, 2, 3, 2, 3, 4, 5, 6, 7, 9, 10, 11, 13, 16, 17, 18, 20, 24, 25, 27, 29, 31, 32)
$oldVal   = @(831, 13652, 13450, 1043, 797, 584, 77, 16, 738, 2125, 66, 79, 103, 363, 495, 154, 1360, 215, 91, 215, 831, 13652, 13450, 1043, 797, 584, 77, 16, 738, 2125, 66, 79, 103, 91, 91, 363, 495, 154, 1360)
$newVal   = @(836, 13665, 13455, 1044, 799, 585, 78, 20, 739, 2128, 71, 80, 106, 364, 496, 155, 1371, 217, 94, 217, 836, 13665, 13455, 1044, 799, 585, 78, 20, 739, 2128, 71, 80, 106, 94, 94, 364, 496, 155, 1371)

$count = 0
$mismatches = 0
for ($i = 0; $i -lt $sheetIdx.Count; $i++) {
    $ws = $wb.Worksheets.Item($sheetIdx[$i])
    $cell = $ws.Cells.Item($rowNum[$i], 6)
    $current = [double]$cell.Value2
    if ($current -ne $oldVal[$i]) {
        $mismatches = $mismatches + 1
    }
    $cell.Value = $newVal[$i]
    $count = $count + 1
}

Write-Output "Done applying $count updates ($mismatches unexpected prior values)."
